$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11167
$ws.Range("E2").Value = 421
$ws.Range("F2").Value = 421
$ws.Range("G2").Value = 291
$ws.Range("H2").Value = 226
$ws.Range("I2").Value = 264
$ws.Range("J2").Value = -38
$ws.Range("K2").Value = 8588
$ws.Range("L2").Value = 4801
$ws.Range("M2").Value = 3787
$ws.Range("N2").Value = 3767
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 303
$ws.Range("R2").Value = -927
$ws.Range("S2").Value = 311
$ws.Range("T2").Value = 781
$ws.Range("U2").Value = -477
$ws.Range("V2").Value = 1372
$ws.Range("W2").Value = 3.77
$ws.Range("X2").Value = 2.02
$ws.Range("Y2").Value = 7.25
$ws.Range("Z2").Value = 2.91
$ws.Range("AA2").Value = 126.79
$ws.Range("AB2").Value = 3668.98
$ws.Range("AC2").Value = 1318
$ws.Range("AD2").Value = 9.71
$ws.Range("AE2").Value = 18924
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.17
$ws.Range("AI2").Value = 11.3
$ws.Range("AJ2").Value = 20054430

# Row 3
$ws.Range("D3").Value = 11991
$ws.Range("E3").Value = 216
$ws.Range("F3").Value = 216
$ws.Range("G3").Value = 247
$ws.Range("H3").Value = 198
$ws.Range("I3").Value = 226
$ws.Range("J3").Value = -28
$ws.Range("K3").Value = 9144
$ws.Range("L3").Value = 5212
$ws.Range("M3").Value = 3932
$ws.Range("N3").Value = 3912
$ws.Range("O3").Value = 21
$ws.Range("P3").Value = 100
$ws.Range("Q3").Value = 287
$ws.Range("R3").Value = -674
$ws.Range("S3").Value = 289
$ws.Range("T3").Value = 1015
$ws.Range("U3").Value = -728
$ws.Range("V3").Value = 1688
$ws.Range("W3").Value = 1.8
$ws.Range("X3").Value = 1.65
$ws.Range("Y3").Value = 5.88
$ws.Range("Z3").Value = 2.23
$ws.Range("AA3").Value = 132.53
$ws.Range("AB3").Value = 3859.56
$ws.Range("AC3").Value = 1125
$ws.Range("AD3").Value = 8.93
$ws.Range("AE3").Value = 19832
$ws.Range("AF3").Value = 0.51
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 2.49
$ws.Range("AI3").Value = 21.85
$ws.Range("AJ3").Value = 20054430

# Row 4
$ws.Range("D4").Value = 11544
$ws.Range("E4").Value = 234
$ws.Range("F4").Value = 234
$ws.Range("G4").Value = 200
$ws.Range("H4").Value = 138
$ws.Range("I4").Value = 144
$ws.Range("J4").Value = -6
$ws.Range("K4").Value = 9525
$ws.Range("L4").Value = 5531
$ws.Range("M4").Value = 3994
$ws.Range("N4").Value = 3971
$ws.Range("O4").Value = 23
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 817
$ws.Range("R4").Value = -703
$ws.Range("S4").Value = 63
$ws.Range("T4").Value = 813
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 1834
$ws.Range("W4").Value = 2.02
$ws.Range("X4").Value = 1.2
$ws.Range("Y4").Value = 3.66
$ws.Range("Z4").Value = 1.48
$ws.Range("AA4").Value = 138.47
$ws.Range("AB4").Value = 3976.7
$ws.Range("AC4").Value = 720
$ws.Range("AD4").Value = 13.16
$ws.Range("AE4").Value = 20355
$ws.Range("AF4").Value = 0.47
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 2.11
$ws.Range("AI4").Value = 27.03
$ws.Range("AJ4").Value = 20054430

# Row 5
$ws.Range("D5").Value = 10359
$ws.Range("E5").Value = -67
$ws.Range("F5").Value = -67
$ws.Range("G5").Value = -69
$ws.Range("H5").Value = -76
$ws.Range("I5").Value = -74
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 9006
$ws.Range("L5").Value = 5197
$ws.Range("M5").Value = 3808
$ws.Range("N5").Value = 3779
$ws.Range("O5").Value = 29
$ws.Range("P5").Value = 100
$ws.Range("Q5").Value = -128
$ws.Range("R5").Value = -490
$ws.Range("S5").Value = 212
$ws.Range("T5").Value = 592
$ws.Range("U5").Value = -720
$ws.Range("V5").Value = 2038
$ws.Range("W5").Value = -0.65
$ws.Range("X5").Value = -0.73
$ws.Range("Y5").Value = -1.92
$ws.Range("Z5").Value = -0.82
$ws.Range("AA5").Value = 136.46
$ws.Range("AB5").Value = 3888.01
$ws.Range("AC5").Value = -371
$ws.Range("AD5").Value = -20.42
$ws.Range("AE5").Value = 19653
$ws.Range("AF5").Value = 0.39
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.32
$ws.Range("AI5").Value = -25.83
$ws.Range("AJ5").Value = 20054430

# Row 6
$ws.Range("D6").Value = 10995
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = -109
$ws.Range("H6").Value = -194
$ws.Range("I6").Value = -192
$ws.Range("K6").Value = 8710
$ws.Range("L6").Value = 5096
$ws.Range("M6").Value = 3614
$ws.Range("N6").Value = 3584
$ws.Range("P6").Value = 102
$ws.Range("Q6").Value = 340
$ws.Range("R6").Value = -404
$ws.Range("S6").Value = 151
$ws.Range("T6").Value = 585
$ws.Range("U6").Value = -245
$ws.Range("V6").Value = 2215
$ws.Range("W6").Value = 0.02
$ws.Range("X6").Value = -1.76
$ws.Range("Y6").Value = -5.21
$ws.Range("Z6").Value = -2.19
$ws.Range("AA6").Value = 141.03
$ws.Range("AB6").Value = 3622.66
$ws.Range("AC6").Value = -955
$ws.Range("AD6").Value = -8.16
$ws.Range("AE6").Value = 18303
$ws.Range("AF6").Value = 0.43
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 0.64
$ws.Range("AI6").Value = -5.11
$ws.Range("AJ6").Value = 20452479

# Clear rows 7-9 data cells (D:AJ), keep A:C
$ws.Range("D7:AJ9").ClearContents()
